# Apply the edit described by the diff: extend the ImportItemAssignedToLocation
# sheet with extra rows (duplicated PCode1 rows, new DRGP partners, and a
# PCode2..PCode10 series), plus update the selected cell in the view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: DRGP-0000002 / N2 / I000038
$ws.Range("A3").Value = "DRGP-0000002"
$ws.Range("B3").Value = "N2"
$ws.Range("C3").Value = "I000038"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

# New row 4: DRGP-0000003 / N3 / I000039 (previously row 3 content: PCode1/ghj/I000042 is gone)
$ws.Range("A4").Value = "DRGP-0000003"
$ws.Range("B4").Value = "N3"
$ws.Range("C4").Value = "I000039"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

# Rows 5-8: PCode1 / ghj / I000055 repeated four times
for ($r = 5; $r -le 8; $r++) {
    $ws.Range("A$r").Value = "PCode1"
    $ws.Range("B$r").Value = "ghj"
    $ws.Range("C$r").Value = "I000055"
    $ws.Range("D$r").Value = 0
    $ws.Range("E$r").Value = 0
}

# Rows 9-17: PCode2..PCode10 / ghj / I000056..I000064
$pcodeNums = 2..10
$itemNums = 56..64
for ($i = 0; $i -lt $pcodeNums.Length; $i++) {
    $r = 9 + $i
    $ws.Range("A$r").Value = "PCode" + $pcodeNums[$i]
    $ws.Range("B$r").Value = "ghj"
    $ws.Range("C$r").Value = "I0000" + $itemNums[$i]
    $ws.Range("D$r").Value = 0
    $ws.Range("E$r").Value = 0
}

# Update the active selection to match the post-edit workbook state.
$ws.Range("I7").Select()
